$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$diffText = "['-   <TipoComunicazione>newwwwww</TipoComunicazione>', '+   <TipoComunicazione>252525</TipoComunicazione>']"

for ($r = 2; $r -le 5; $r++) {
    # VALORE DA MODIFICARE column (C) was empty; now filled with "newwwwww"
    $ws.Cells.Item($r, 3).Value = "newwwwww"

    # VALORE MODIFICATO column (D) changes from "99999" to "252525" (kept as text)
    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = "252525"

    # DIFFERENZE column (F) gets the new, clearer diff-style text
    $ws.Cells.Item($r, 6).Value = $diffText
}
